$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Map"
$ws.Range("D2").Value = "Map"
$ws.Range("F2").Value = "'0"
$ws.Range("G2").Value = "'0"
$ws.Range("E2").Value = "Binder"
$ws.Range("H2").Value = "Map"
$ws.Range("K2").Value = "'"
